$d = $word.ActiveDocument

$d.Content.Find.Execute('2023-08-23 Wednesday', $true, $false, $false, $false, $false, $true, 1, $false, '2023-08-24 Thursday', 2) | Out-Null
$d.Content.Find.Execute('33×67=', $true, $false, $false, $false, $false, $true, 1, $false, '11×50=', 2) | Out-Null
$d.Content.Find.Execute('36×93=', $true, $false, $false, $false, $false, $true, 1, $false, '18×60=', 2) | Out-Null
$d.Content.Find.Execute('97×11=', $true, $false, $false, $false, $false, $true, 1, $false, '49×87=', 2) | Out-Null
$d.Content.Find.Execute('46×24=', $true, $false, $false, $false, $false, $true, 1, $false, '98×51=', 2) | Out-Null
$d.Content.Find.Execute('35×42=', $true, $false, $false, $false, $false, $true, 1, $false, '38×19=', 2) | Out-Null
$d.Content.Find.Execute('15×15=', $true, $false, $false, $false, $false, $true, 1, $false, '57×27=', 2) | Out-Null
$d.Content.Find.Execute('97×26=', $true, $false, $false, $false, $false, $true, 1, $false, '20×44=', 2) | Out-Null
$d.Content.Find.Execute('71×37=', $true, $false, $false, $false, $false, $true, 1, $false, '65×21=', 2) | Out-Null
$d.Content.Find.Execute('63×19=', $true, $false, $false, $false, $false, $true, 1, $false, '43×27=', 2) | Out-Null
$d.Content.Find.Execute('32×26=', $true, $false, $false, $false, $false, $true, 1, $false, '76×66=', 2) | Out-Null
$d.Content.Find.Execute('83×47=', $true, $false, $false, $false, $false, $true, 1, $false, '67×59=', 2) | Out-Null
$d.Content.Find.Execute('42×29=', $true, $false, $false, $false, $false, $true, 1, $false, '61×24=', 2) | Out-Null
$d.Content.Find.Execute('87×85=', $true, $false, $false, $false, $false, $true, 1, $false, '70×99=', 2) | Out-Null
$d.Content.Find.Execute('82×50=', $true, $false, $false, $false, $false, $true, 1, $false, '55×71=', 2) | Out-Null
$d.Content.Find.Execute('45×81=', $true, $false, $false, $false, $false, $true, 1, $false, '38×37=', 2) | Out-Null
$d.Content.Find.Execute('31×24=', $true, $false, $false, $false, $false, $true, 1, $false, '73×86=', 2) | Out-Null
$d.Content.Find.Execute('60×11=', $true, $false, $false, $false, $false, $true, 1, $false, '21×15=', 2) | Out-Null
$d.Content.Find.Execute('56×72=', $true, $false, $false, $false, $false, $true, 1, $false, '81×60=', 2) | Out-Null
$d.Content.Find.Execute('35×38=', $true, $false, $false, $false, $false, $true, 1, $false, '41×73=', 2) | Out-Null
$d.Content.Find.Execute('12×38=', $true, $false, $false, $false, $false, $true, 1, $false, '79×53=', 2) | Out-Null
$d.Content.Find.Execute('48×50=', $true, $false, $false, $false, $false, $true, 1, $false, '79×21=', 2) | Out-Null
$d.Content.Find.Execute('40×39=', $true, $false, $false, $false, $false, $true, 1, $false, '37×46=', 2) | Out-Null
$d.Content.Find.Execute('86×15=', $true, $false, $false, $false, $false, $true, 1, $false, '24×20=', 2) | Out-Null
$d.Content.Find.Execute('13×20=', $true, $false, $false, $false, $false, $true, 1, $false, '67×59=', 2) | Out-Null
$d.Content.Find.Execute('63×37=', $true, $false, $false, $false, $false, $true, 1, $false, '43×68=', 2) | Out-Null
